# Generate Report for Handoff
# b.md's handback xliff arrived but does not match the latest source
# version, so its status flips from "Handed back: in sync with en-US"
# to "Ready for handoff" on the Overview sheet and on each language
# sheet (zh-cn / de-de); the corresponding "latest handback" metadata
# for b.md is updated accordingly.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/74e387e32b5ac13be3d15a55947823cd8ad489cf/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4abadff3bf3c8b2e6cb6a30c9d2173fb2e51639c/e2e/b.md."

# ---- Overview sheet: b.md row (row 3) ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("G3").Value = "2016-08-21 20:49:18"

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
# Leading apostrophe keeps "False" stored as text (not auto-coerced to a
# boolean); resetting the style back to Normal drops the quote-prefix
# formatting flag that the apostrophe otherwise leaves behind.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-21 20:49:13"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-21 20:49:18"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.17
